# This script edits noh_divorce.docx:
#  - Block 1 (Plaintiff name cell): replaces the
#      {% if user_ask_role == "plaintiff" %}{{ users[0].name_full() | upper }}{% else %}{{ other_parties[0].name_full() | upper }}{% endif %}
#    jinja2 conditional with a simple
#      {{ plaintiffs[0].name_full() | upper }}
#    expression.
#  - Block 2 (Defendant name cell): replaces the
#      {% if user_ask_role == "plaintiff" %}{{ other_parties[0].name_full() | upper }}{% else %}{{ users[0].name_full() | upper }}{% endif %}
#    jinja2 conditional with a simple
#      {{ defendants[0].name_full() | upper }}
#    expression (keeping the trailing "/" and adding a space before it).

$d = $word.ActiveDocument

$LDQ = [char]8220   # left double quotation mark "
$RDQ = [char]8221   # right double quotation mark "

# ---------------------------------------------------------------------
# Block 1: the Plaintiff-side name cell
# ---------------------------------------------------------------------

# 1a. Remove the leading "{% if user_ask_role == "plaintiff" %}" tag
#     (this occurs twice in the doc; this is the first occurrence).
$ifTag = "{% if user_ask_role == " + $LDQ + "plaintiff" + $RDQ + " %}"
$rng = $d.Content
$rng.Find.Execute($ifTag, $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

# 1b. Turn "users" (inside the bold "{{ users[0].name_full() ..." run)
#     into "plaintiffs".
$rng = $d.Content
$rng.Find.Execute("{{ users[0].name_full()", $true, $false, $false, $false, $false, $true, 1, $false, "{{ plaintiffs[0].name_full()", 1) | Out-Null

# 1c. Remove the "{% else %}{{ other_parties[0].name_full() | upper }}{% endif %}" tail.
$rng = $d.Content
$rng.Find.Execute("{% else %}{{ other_parties[0].name_full() | upper }}{% endif %}", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

# ---------------------------------------------------------------------
# Block 2: the Defendant-side name cell
# ---------------------------------------------------------------------

# 2a. Remove the leading "{% if user_ask_role == "plaintiff" %}" tag
#     (second/remaining occurrence).
$rng = $d.Content
$rng.Find.Execute($ifTag, $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

# 2b. Turn "other_parties" (inside the bold "{{ other_parties[0].name_full() ..." run)
#     into "defendants".
$rng = $d.Content
$rng.Find.Execute("other_parties", $true, $false, $false, $false, $false, $true, 1, $false, "defendants", 1) | Out-Null

# 2c. Remove the "{% else %}{{ users[0].name_full() | upper }}" portion, and turn the
#     trailing "{% endif %}" into a single space (keeping its run/formatting alive),
#     so the cell reads "...}} /" instead of "...}}{% endif %}/".
$rng = $d.Content
$rng.Find.Execute("{% else %}{{ users[0].name_full() | upper }}{% endif %}", $true, $false, $false, $false, $false, $true, 1, $false, " ", 1) | Out-Null
